$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# D14: flip sign to negative
$ws.Range("D14").Value = -45752811.060000002

# D16: flip sign to negative
$ws.Range("D16").Value = -50601311.960000001

# D18: replace static value with SUM formula
$ws.Range("D18").Formula = "=SUM(D12:D17)"

# D19: flip sign to negative
$ws.Range("D19").Value = -383100000

# D21: replace static value with SUM formula
$ws.Range("D21").Formula = "=SUM(D18:D20)"

# D22: replace static value
$ws.Range("D22").Value = -322025175

$wb.Save()
